# "first multi weight version"
#
# 1. Rename "ComputeWeighting" -> "CheckWeighting"
# 2. Insert a new sheet "default_Jnd" right after "CheckWeighting" (before "JndPer")
#    and populate it with Abs/Per/value/weight columns.
# 3. Tidy up a few sheet view / selection states and the CheckWeighting page setup.

$wb = $excel.ActiveWorkbook

# --- 1. Rename ComputeWeighting -> CheckWeighting -------------------------
$checkWs = $wb.Worksheets.Item("ComputeWeighting")
$checkWs.Name = "CheckWeighting"

# Refresh its view: zoom to 70% and select E3:X57 instead of the old N30.
$checkWs.Activate() | Out-Null
$win = $wb.Windows.Item(1)
$win.Zoom = 70
$checkWs.Range("E3:X57").Select() | Out-Null

# Touching PageSetup drops the now-stale printer-settings relationship on
# this sheet (keeps the orientation the same, just forces a rewrite).
$checkWs.PageSetup.Orientation = 1

# --- 2. Insert new "default_Jnd" sheet between CheckWeighting and JndPer --
$jndWs = $wb.Worksheets.Add($null, $checkWs)
$jndWs.Name = "default_Jnd"

$jndWs.Range("B1").Value = "Abs"
$jndWs.Range("C1").Value = "Per"
$jndWs.Range("D1").Value = "value"
$jndWs.Range("E1").Value = "weight"

$jndWs.Range("A2").Value = "Fare"
$jndWs.Range("B2").Value = 15
$jndWs.Range("C2").Value = 20
$jndWs.Range("D2").Value = 0.2
$jndWs.Range("E2").Value = 1

$jndWs.Range("A3").Value = "Travel"
$jndWs.Range("B3").Value = 120
$jndWs.Range("C3").Value = 100
$jndWs.Range("D3").Value = 0.2
$jndWs.Range("E3").Value = 1

$jndWs.Range("A4").Value = "Wait"
$jndWs.Range("B4").Value = 100
$jndWs.Range("C4").Value = 100
$jndWs.Range("D4").Value = 0
$jndWs.Range("E4").Value = 1

$jndWs.Range("A5").Value = "Transfer"
$jndWs.Range("B5").Value = 10
$jndWs.Range("C5").Value = 10
$jndWs.Range("D5").Value = 1
$jndWs.Range("E5").Value = 1

$jndWs.Range("A6").Value = "Walk"
$jndWs.Range("B6").Value = -1
$jndWs.Range("C6").Value = -3
$jndWs.Range("D6").Value = 1
$jndWs.Range("E6").Value = 1

# --- 3. Misc selection tidy-ups on other sheets ---------------------------
$paraWs = $wb.Worksheets.Item("Para")
$paraWs.Activate() | Out-Null
$paraWs.Range("B5").Select() | Out-Null

$pasWs = $wb.Worksheets.Item("PasWeight")
$pasWs.Activate() | Out-Null
$pasWs.Range("B2:B6").Select() | Out-Null

$jndAbsWs = $wb.Worksheets.Item("JndAbs")
$jndAbsWs.Activate() | Out-Null
$jndAbsWs.Range("A1:C6").Select() | Out-Null

# Leave "default_Jnd" as the active sheet/tab, with I15 selected.
$jndWs.Activate() | Out-Null
$jndWs.Range("I15").Select() | Out-Null
